$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each cell's existing "a÷b=" quiz value is replaced with a new one.
# Editing is done per-cell via the Cell's own Range so that the two
# duplicated source strings (e.g. "42÷2=") and the accidental
# old-value/new-value collisions between cells (e.g. "35÷3=" -> "53÷4="
# later followed by another cell's "53÷4=" -> "76÷8=") each land on the
# correct cell instead of a document-wide find/replace clobbering them.
$replacements = @(
    @{ Row = 1;  Col = 1; Old = "30÷2="; New = "20÷7=" },
    @{ Row = 1;  Col = 2; Old = "85÷9="; New = "95÷8=" },
    @{ Row = 1;  Col = 3; Old = "42÷2="; New = "46÷4=" },
    @{ Row = 1;  Col = 4; Old = "98÷2="; New = "81÷7=" },
    @{ Row = 1;  Col = 5; Old = "95÷2="; New = "68÷6=" },

    @{ Row = 5;  Col = 1; Old = "10÷5="; New = "87÷9=" },
    @{ Row = 5;  Col = 2; Old = "41÷2="; New = "48÷2=" },
    @{ Row = 5;  Col = 3; Old = "50÷9="; New = "55÷5=" },
    @{ Row = 5;  Col = 4; Old = "65÷6="; New = "69÷2=" },
    @{ Row = 5;  Col = 5; Old = "76÷3="; New = "33÷4=" },

    @{ Row = 9;  Col = 1; Old = "30÷4="; New = "13÷2=" },
    @{ Row = 9;  Col = 2; Old = "35÷3="; New = "53÷4=" },
    @{ Row = 9;  Col = 3; Old = "88÷6="; New = "25÷4=" },
    @{ Row = 9;  Col = 4; Old = "71÷8="; New = "35÷9=" },
    @{ Row = 9;  Col = 5; Old = "56÷4="; New = "54÷5=" },

    @{ Row = 13; Col = 1; Old = "31÷5="; New = "18÷5=" },
    @{ Row = 13; Col = 2; Old = "71÷6="; New = "30÷3=" },
    @{ Row = 13; Col = 3; Old = "54÷4="; New = "70÷4=" },
    @{ Row = 13; Col = 4; Old = "53÷5="; New = "20÷3=" },
    @{ Row = 13; Col = 5; Old = "87÷2="; New = "39÷6=" },

    @{ Row = 17; Col = 1; Old = "63÷9="; New = "78÷7=" },
    @{ Row = 17; Col = 2; Old = "85÷6="; New = "84÷5=" },
    @{ Row = 17; Col = 3; Old = "42÷2="; New = "13÷6=" },
    @{ Row = 17; Col = 4; Old = "53÷4="; New = "76÷8=" },
    @{ Row = 17; Col = 5; Old = "18÷5="; New = "19÷8=" }
)

foreach ($r in $replacements) {
    $cell = $t.Cell($r.Row, $r.Col)
    $cell.Range.Text = $r.New
}
